$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.763.39"
$ws.Range("E2").Value = "  -0.41%  "
$ws.Range("D3").Value = "1.596.20"
$ws.Range("E3").Value = "  -1.65%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "209.25"
$ws.Range("E5").Value = "  -1.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.502"
$ws.Range("E6").Value = "  -2.00%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.39"
$ws.Range("E8").Value = "  -2.48%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.253"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0593"
$ws.Range("E10").Value = "  -1.65%  "
$ws.Range("E11").Value = "  -1.57%  "
$ws.Range("D12").Value = "1.823.34"
$ws.Range("E12").Value = "  -1.62%  "
$ws.Range("D13").Value = "1.602.05"
$ws.Range("E13").Value = "  -0.95%  "
$ws.Range("E14").Value = "  -2.22%  "
$ws.Range("E15").Value = "  -3.32%  "
$ws.Range("D16").Value = "27.756.51"
$ws.Range("E16").Value = "  -0.39%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.51"
$ws.Range("E17").Value = "  -1.56%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "219.94"
$ws.Range("E18").Value = "  -2.96%  "
$ws.Range("E19").Value = "  -2.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.37"
$ws.Range("E20").Value = "  -2.79%  "
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("E22").Value = "  -3.43%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.83"
$ws.Range("E23").Value = "  -0.96%  "
$ws.Range("E24").Value = "  -3.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.05"
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.18"
$ws.Range("E26").Value = "  +4.17%  "
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.19"
$ws.Range("E28").Value = "  -0.94%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.106"
$ws.Range("E29").Value = "  -3.75%  "
$ws.Range("E30").Value = "  -0.84%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0473"
$ws.Range("E31").Value = "  -1.24%  "
$ws.Range("E32").Value = "  -4.06%  "
$ws.Range("D33").Value = "1.376.92"
$ws.Range("E33").Value = "  -3.02%  "
$ws.Range("E34").Value = "  -2.89%  "
$ws.Range("E35").Value = "  -3.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.976"
$ws.Range("E36").Value = "  -0.40%  "
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("E38").Value = "  -0.50%  "
$ws.Range("E39").Value = "  -2.77%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.828"
$ws.Range("E40").Value = "  -1.83%  "
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.975"
$ws.Range("E42").Value = "  -2.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "64.62"
$ws.Range("E43").Value = "  -0.65%  "
$ws.Range("E44").Value = "  +2.49%  "
$ws.Range("E45").Value = "  -1.31%  "
$ws.Range("E46").Value = "  -1.87%  "
$ws.Range("D47").Value = "1.734.71"
$ws.Range("E47").Value = "  -1.68%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.79"
$ws.Range("E48").Value = "  -2.80%  "
$ws.Range("E49").Value = "  -0.65%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0969"
$ws.Range("E50").Value = "  -2.72%  "
$ws.Range("E51").Value = "  -0.75%  "
